$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($sheet, $addr, $val) {
    $c = $sheet.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Row 2 ---
Set-CellText $ws "D2" "41.202.76"
$ws.Range("E2").Value = "  -6.47%  "

# --- Row 3 ---
Set-CellText $ws "D3" "2.212.18"
$ws.Range("E3").Value = "  -6.77%  "

# --- Row 4 ---
$ws.Range("E4").Value = "  +0.00%  "

# --- Row 5 ---
Set-CellText $ws "D5" "243.05"
$ws.Range("E5").Value = "  +1.11%  "

# --- Row 6 ---
Set-CellText $ws "D6" "0.625"
$ws.Range("E6").Value = "  -7.30%  "

# --- Row 7 ---
Set-CellText $ws "D7" "70.16"
$ws.Range("E7").Value = "  -5.81%  "

# --- Row 8 ---
$ws.Range("E8").Value = "  -0.08%  "

# --- Row 9 ---
Set-CellText $ws "D9" "0.550"
$ws.Range("E9").Value = "  -9.79%  "

# --- Row 10 ---
Set-CellText $ws "D10" "37.29"
$ws.Range("E10").Value = "  +0.23%  "

# --- Row 11 ---
Set-CellText $ws "D11" "0.0953"
$ws.Range("E11").Value = "  -7.27%  "

# --- Row 12 ---
Set-CellText $ws "D12" "57.74"
$ws.Range("E12").Value = "  -3.73%  "

# --- Row 13 ---
$ws.Range("E13").Value = "  -4.34%  "

# --- Row 14 ---
Set-CellText $ws "D14" "6.68"
$ws.Range("E14").Value = "  -8.74%  "

# --- Row 15 ---
Set-CellText $ws "D15" "2.537.38"
$ws.Range("E15").Value = "  -7.04%  "

# --- Row 16 ---
Set-CellText $ws "D16" "14.74"
$ws.Range("E16").Value = "  -10.25%  "

# --- Row 17 ---
Set-CellText $ws "D17" "0.838"
$ws.Range("E17").Value = "  -9.84%  "

# --- Row 18 ---
Set-CellText $ws "D18" "2.202.97"
$ws.Range("E18").Value = "  -7.38%  "

# --- Row 19 ---
Set-CellText $ws "D19" "41.103.60"
$ws.Range("E19").Value = "  -6.61%  "

# --- Row 20 ---
Set-CellText $ws "D20" "0.0₃0947"
$ws.Range("E20").Value = "  -8.52%  "

# --- Row 21 ---
Set-CellText $ws "D21" "72.65"
$ws.Range("E21").Value = "  -6.92%  "

# --- Row 22 ---
Set-CellText $ws "D22" "6.07"
$ws.Range("E22").Value = "  -8.36%  "

# --- Row 23 ---
Set-CellText $ws "D23" "231.14"
$ws.Range("E23").Value = "  -9.23%  "

# --- Row 24 ---
Set-CellText $ws "D24" "2.02"
$ws.Range("E24").Value = "  +7.01%  "

# --- Row 25 ---
Set-CellText $ws "D25" "1.00"
$ws.Range("E25").Value = "  +0.03%  "

# --- Row 26 ---
Set-CellText $ws "D26" "3.59"
$ws.Range("E26").Value = "  -5.10%  "

# --- Row 27 ---
$ws.Range("E27").Value = "  -3.58%  "

# --- Row 28 ---
Set-CellText $ws "D28" "2.18"
$ws.Range("E28").Value = "  -5.04%  "

# --- Row 29 ---
Set-CellText $ws "D29" "9.74"
$ws.Range("E29").Value = "  -8.08%  "

# --- Row 30 ---
Set-CellText $ws "D30" "171.21"
$ws.Range("E30").Value = "  -2.60%  "

# --- Row 31 ---
Set-CellText $ws "D31" "20.39"
$ws.Range("E31").Value = "  -9.09%  "

# --- Row 32 ---
Set-CellText $ws "D32" "0.119"
$ws.Range("E32").Value = "  -8.93%  "

# --- Row 33 ---
$ws.Range("E33").Value = "  -7.91%  "

# --- Row 34 ---
Set-CellText $ws "D34" "0.0705"
$ws.Range("E34").Value = "  -7.15%  "

# --- Row 35 ---
Set-CellText $ws "D35" "5.16"
$ws.Range("E35").Value = "  -5.03%  "

# --- Row 36 ---
$ws.Range("E36").Value = "  -9.98%  "

# --- Row 37 ---
Set-CellText $ws "D37" "3.86"
$ws.Range("E37").Value = "  +1.35%  "

# --- Row 38 ---
Set-CellText $ws "D38" "23.86"
$ws.Range("E38").Value = "  +15.07%  "

# --- Row 39 (was LidoDAOToken -> now VeChain) ---
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText $ws "D39" "0.0275"
$ws.Range("E39").Value = "  -2.92%  "

# --- Row 40 (was VeChain -> now LidoDAOToken) ---
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-CellText $ws "D40" "2.27"
$ws.Range("E40").Value = "  -5.85%  "

# --- Row 41 ---
Set-CellText $ws "D41" "5.82"
$ws.Range("E41").Value = "  -12.64%  "

# --- Row 42 ---
Set-CellText $ws "D42" "64.09"
$ws.Range("E42").Value = "  -2.78%  "

# --- Row 43 ---
Set-CellText $ws "D43" "4.90"
$ws.Range("E43").Value = "  -11.92%  "

# --- Row 44 ---
Set-CellText $ws "D44" "0.195"
$ws.Range("E44").Value = "  -4.63%  "

# --- Row 45 ---
Set-CellText $ws "D45" "8.59"
$ws.Range("E45").Value = "  -5.70%  "

# --- Row 46 (was BinanceUSD -> now Cronos) ---
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText $ws "D46" "0.0998"
$ws.Range("E46").Value = "  -7.63%  "

# --- Row 47 (was Cronos -> now BinanceUSD) ---
$ws.Range("B47").Value = "BinanceUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-CellText $ws "D47" "1.00"
$ws.Range("E47").Value = "  +0.07%  "

# --- Row 48 ---
Set-CellText $ws "D48" "10.56"
$ws.Range("E48").Value = "  +8.95%  "

# --- Row 49 ---
$ws.Range("E49").Value = "  -0.08%  "

# --- Row 50 ---
Set-CellText $ws "D50" "1.17"
$ws.Range("E50").Value = "  -6.55%  "

# --- Row 51 ---
$ws.Range("E51").Value = "  -6.31%  "
